$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.406.46'
$ws.Range("E2").Value = '  +0.73%  '
$ws.Range("D3").Value = '1.625.27'
$ws.Range("E3").Value = '  +1.28%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''212.64'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("D9").Value = '''0.0617'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.47%  '
$ws.Range("D10").Value = '''18.91'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.08%  '
$ws.Range("D11").Value = '''0.0834'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.11%  '
$ws.Range("D12").Value = '1.852.22'
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("D13").Value = '1.625.80'
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("E14").Value = '  +0.76%  '
$ws.Range("D15").Value = '''0.521'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.07%  '
$ws.Range("D16").Value = '26.395.38'
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("D17").Value = '''62.66'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.20%  '
$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").Value = '''202.94'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.46%  '
$ws.Range("D21").Value = '''4.27'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").Value = '''9.34'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.79%  '
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("E24").Value = '  -3.50%  '
$ws.Range("D25").Value = '''144.63'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  -2.86%  '
$ws.Range("D28").Value = '''15.21'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("E29").Value = '  +0.85%  '
$ws.Range("E30").Value = '  +5.54%  '
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("E32").Value = '  +1.52%  '
$ws.Range("D33").Value = '''2.93'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("D34").Value = '''1.50'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.08%  '
$ws.Range("E35").Value = '  +2.28%  '
$ws.Range("D36").Value = '1.158.59'
$ws.Range("E36").Value = '  +1.86%  '
$ws.Range("E37").Value = '  +0.63%  '
$ws.Range("D38").Value = '''0.805'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.60%  '
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").Value = '''0.497'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("E42").Value = '  +3.57%  '
$ws.Range("E43").Value = '  -0.28%  '
$ws.Range("D44").Value = '1.763.55'
$ws.Range("E44").Value = '  +1.39%  '
$ws.Range("D45").Value = '''92.05'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = '''1.52'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.39%  '
$ws.Range("E47").Value = '  +9.04%  '
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("D51").Value = '''0.998'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.42%  '
